$d = $word.ActiveDocument

# 1) Fix typo "mixs" -> "mixes" in "It mixs the data smooth."
$d.Content.Find.Execute("It mixs the data smooth.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "It mixes the data smooth.", 2) | Out-Null

# 2) Split "log_Budget <- log10(mydata$Budget) log_Budget" (FirstParagraph style)
#    into two separate paragraphs:
#      "log_Budget <- log10(mydata$Budget) "
#      "log_Budget"
$d.Content.Find.Execute("log_Budget <- log10(mydata`$Budget) log_Budget", $false, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "log_Budget <- log10(mydata`$Budget) ^plog_Budget", 2) | Out-Null

# 3) Insert a blank paragraph before, and split
#    "log_Budget <- log10(mydata$Audience.Ratings) log_Budget" (BodyText style)
#    into two separate paragraphs:
#      "log_Budget <- log10(mydata$Audience.Ratings) "
#      "log_Budget"
$d.Content.Find.Execute("log_Budget <- log10(mydata`$Audience.Ratings) log_Budget", $false, $false, `
                         $false, $false, $false, $true, 1, $false, `
                         "^plog_Budget <- log10(mydata`$Audience.Ratings) ^plog_Budget", 2) | Out-Null
